# Trying outlier detection adding ao only
# Deletes the LS-outlier rows (LS2001Sep, LS2017Aug, LS2017Oct, LS2020Mar,
# LS2020Jun, LS2020Aug, LS2020Dec, LS2021Mar) and keeps only the AO (additive
# outlier) rows plus the ARIMA model terms, refreshing the coefficient table
# with the new model's summary statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (LS2001Sep) entirely - AO2008Sep (currently row 3)
# becomes the new row 2 along with everything below it.
$ws.Range("A2:G2").EntireRow.Delete()

# After the shift above, the rows (by original label) now sitting at:
#   2: AO2008Sep   3: LS2017Aug  4: LS2017Oct  5: AO2020Mar  6: LS2020Mar
#   7: LS2020Jun   8: LS2020Aug  9: LS2020Dec  10: LS2021Mar 11: ar.L1
#   12: ma.L1      13: ma.L2     14: ar.S.L12  15: ma.S.L12  16: sigma2
# Remove the remaining LS* rows (now at 3,4,6,7,8,9,10), keeping AO2008Sep (2)
# and AO2020Mar (which is row 5 before these deletions).
$ws.Range("A3:G4").EntireRow.Delete()
# Now: 2 AO2008Sep, 3 AO2020Mar, 4 LS2020Mar, 5 LS2020Jun, 6 LS2020Aug,
#      7 LS2020Dec, 8 LS2021Mar, 9 ar.L1, 10 ma.L1, 11 ma.L2, 12 ar.S.L12,
#      13 ma.S.L12, 14 sigma2
$ws.Range("A4:G8").EntireRow.Delete()
# Now: 2 AO2008Sep, 3 AO2020Mar, 4 ar.L1, 5 ma.L1, 6 ma.L2, 7 ar.S.L12,
#      8 ma.S.L12, 9 sigma2

# Refresh all of the remaining data rows with the new model's statistics.
$ws.Range("A2").Value = "AO2008Sep"
$ws.Range("B2").Value = -134900
$ws.Range("C2").Value = 38900
$ws.Range("D2").Value = -3.471
$ws.Range("E2").Value = 0.001
$ws.Range("F2").Value = -211000
$ws.Range("G2").Value = -58700

$ws.Range("A3").Value = "AO2020Mar"
$ws.Range("B3").Value = 22960
$ws.Range("C3").Value = 10500
$ws.Range("D3").Value = 2.192
$ws.Range("E3").Value = 0.028
$ws.Range("F3").Value = 2429.908
$ws.Range("G3").Value = 43500

$ws.Range("A4").Value = "ar.L1"
$ws.Range("B4").Value = 0.8799
$ws.Range("C4").Value = 0.055
$ws.Range("D4").Value = 15.995
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.772
$ws.Range("G4").Value = 0.988

$ws.Range("A5").Value = "ma.L1"
$ws.Range("B5").Value = -0.8314
$ws.Range("C5").Value = 0.079
$ws.Range("D5").Value = -10.475
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = -0.987
$ws.Range("G5").Value = -0.676

$ws.Range("A6").Value = "ma.L2"
$ws.Range("B6").Value = -0.1528
$ws.Range("C6").Value = 0.063
$ws.Range("D6").Value = -2.43
$ws.Range("E6").Value = 0.015
$ws.Range("F6").Value = -0.276
$ws.Range("G6").Value = -0.03

$ws.Range("A7").Value = "ar.S.L12"
$ws.Range("B7").Value = 0.9617
$ws.Range("C7").Value = 0.027
$ws.Range("D7").Value = 36.099
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0.909
$ws.Range("G7").Value = 1.014

$ws.Range("A8").Value = "ma.S.L12"
$ws.Range("B8").Value = -0.8188
$ws.Range("C8").Value = 0.064
$ws.Range("D8").Value = -12.871
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = -0.943
$ws.Range("G8").Value = -0.694

$ws.Range("A9").Value = "sigma2"
$ws.Range("B9").Value = 2361000000
$ws.Range("C9").Value = 0.731
$ws.Range("D9").Value = 3230000000
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 2360000000
$ws.Range("G9").Value = 2360000000
